$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case identifiers (column A) for the two new rows
$ws.Range("A5").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields"
$ws.Range("A6").Value = "103_TruckInsurance_003_InsurantData_002_FieldHintsAndErrors"

# Row 5 detail columns
$ws.Range("D5").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields_FillFirstName"
$ws.Range("E5").Value = "103_TruckInsurance_003_InsurantData_001_MandatoryFields_CheckFilledFirstName"

# Column B repeated value for the new rows
$ws.Range("B5:B7").Value = "Button Next from Page VehicleData"

# Row 5 / row 6 C column
$ws.Range("C5").Value = "Insurant Data Page check for open mandatory fields"
$ws.Range("C6").Value = "Insurant Data Page check for hints regarding mandatory fields"

# Row 6 remaining columns
$ws.Range("D6").Value = "103_TruckInsurance_003_InsurantData_002_EnterValuesInWrongFormat"
$ws.Range("E6").Value = "Insurant Data Page check error hint formatting"
$ws.Range("F6").Value = "103_TruckInsurance_003_InsurantData_002_EnterValuesInWrongFormat Part 2"
$ws.Range("G6").Value = "Insurant Data Page check error hint formatting Part 2"

# Row 7 - new list-contents test case
$ws.Range("A7").Value = "103_TruckInsurance_003_InsurantData_003_ListContents"
$ws.Range("C7").Value = "103_TruckInsurance_003_InsurantData_003_ListContents"

# Re-fit column widths to match content (mirrors original authoring workflow)
$ws.Range("D1").ColumnWidth = 70.333333333
$ws.Range("E1").ColumnWidth = 70.333333333
$ws.Range("G1").ColumnWidth = 64.833333333

# Update selection to match final saved state
$ws.Range("C7").Select()
